# Daily attendance processing - 2026-01-01 21:55:27
# Normalize "Recorded By" (column G) entries so the "System" token's
# position matches the canonical ordering used across the report.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Exact-value swap map: any G cell whose text matches a key is rewritten
# to the corresponding value.
$map = @{
    "System, backup@backdoor.com"         = "backup@backdoor.com, System"
    "System, backup@backdoor.com, system" = "backup@backdoor.com, System, system"
    "dnasr281@gmail.com, System"          = "System, dnasr281@gmail.com"
}

$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2
    if ($null -ne $val -and $map.ContainsKey($val)) {
        $cell.Value = $map[$val]
    }
}
